$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-obsolete forecast values in C2 and C3 (removed entirely in the diff)
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Updated (corrected) forecast values after the naive component forecaster bug fix
$ws.Range("C4").Value = 0.348613976222456
$ws.Range("C5").Value = -0.1384957661262676
$ws.Range("E5").Value = 0.5922117994852982

$ws.Range("C6").Value = 1.566479473280191
$ws.Range("E6").Value = 0.5295895589954469

$ws.Range("C7").Value = 0.7307568962937161

$ws.Range("C8").Value = 0.8188188121642126

$ws.Range("E9").Value = 0.9697679806506043

$ws.Range("C10").Value = 1.9846842782967

$ws.Range("E11").Value = 1.083941060573257

$ws.Range("C13").Value = 1.064321453542272
$ws.Range("E13").Value = 0.8791375467670504

$ws.Range("C14").Value = 1.361817904277718
$ws.Range("E14").Value = 1.226035857429419

$ws.Range("C15").Value = -4.352425014431327
$ws.Range("E15").Value = -1.222056059052357

$ws.Range("E16").Value = 0.6980411378030515

$ws.Range("E17").Value = 1.605691900741091

$ws.Range("C18").Value = -0.9008525709169657

$ws.Range("C19").Value = 0.2738544794132602
$ws.Range("E19").Value = 0.6247505135484221
